$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.157695889472961
$ws.Range("B1").Value = 2.379462003707886
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.393099069595337
$ws.Range("E1").Value = 1.222819924354553
